$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update order dates in column B (rows 4-13): shift from June 2010 to
# consecutive days in June 2018.
$dates = @(43253, 43254, 43255, 43256, 43257, 43258, 43259, 43260, 43261, 43262)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}

# Fix a typo ("usasemos" -> "usásemos") in the sacks question.
$ws.Range("A26").Value = "Si para servir los kilos de nuestros pedidos usásemos sacos de 25 Kg, cuántos sacos harían falta para servir los pedidos del día 14?"

# The "Importe de los pedidos anteriores" label now mentions the new cutoff
# date.
$ws.Range("A19").Value = "Importe de los pedidos anteriores al 8/6/2018"

# The "¿Cuántos pedidos se realizaron antes del 8/6/2010?" question is
# removed (row 18 left blank).
$ws.Range("A18").Value = ""

# Leave the selection where the author left it when saving.
$ws.Range("K34").Select()
